$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 232
$ws.Range("I4").Value = 232
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 232
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -118
$ws.Range("N4").ClearContents()
$ws.Range("H8").Value = 16
$ws.Range("I8").Value = 16
$ws.Range("K8").Value = 48
$ws.Range("M8").Value = 91
$ws.Range("H9").Value = 28.428572
$ws.Range("I9").Value = 28.166666
$ws.Range("K9").Value = 28.166666
$ws.Range("M9").Value = 140.833334
$ws.Range("H15").Value = 451.4
$ws.Range("I15").Value = 451.4
$ws.Range("K15").Value = 1354.2
$ws.Range("M15").Value = -1185.2
$ws.Range("H31").Value = 1055.5
$ws.Range("I31").Value = 1055.5
$ws.Range("K31").Value = 3166.5
$ws.Range("M31").Value = -2936.5
$ws.Range("H116").Value = 2999.3333
$ws.Range("I116").Value = 2499.5
$ws.Range("K116").Value = 2499.5
$ws.Range("M116").Value = 942.5
$ws.Range("H125").Value = 1066
$ws.Range("I125").Value = 632
$ws.Range("J125").Value = 1500
$ws.Range("K125").Value = 5688
$ws.Range("L125").Value = 13500
$ws.Range("M125").Value = -3228
$ws.Range("N125").Value = -18420
$ws.Range("H138").Value = 5199.6
$ws.Range("I138").Value = 2000
$ws.Range("J138").Value = 5999.5
$ws.Range("K138").Value = 6000
$ws.Range("L138").Value = 17998.5
$ws.Range("M138").Value = -860
$ws.Range("N138").Value = -28278.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4500
$ws.Range("I2").Value = 4500
$ws.Range("K2").Value = 4500
$ws.Range("M2").Value = -4387
$ws.Range("H5").Value = 27.5
$ws.Range("I5").Value = 25
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 25
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 87
$ws.Range("N5").Value = -254
$ws.Range("H8").Value = 3333
$ws.Range("J8").Value = 3333
$ws.Range("L8").Value = 3333
$ws.Range("N8").Value = -3621
$ws.Range("H10").Value = 4999
$ws.Range("J10").Value = 5999
$ws.Range("L10").Value = 5999
$ws.Range("N10").Value = -6339
$ws.Range("H11").Value = 4999
$ws.Range("J11").Value = 4999
$ws.Range("L11").Value = 4999
$ws.Range("N11").Value = -5287
$ws.Range("H19").Value = 498599.5
$ws.Range("I19").Value = 990949.5
$ws.Range("J19").Value = 6249.5
$ws.Range("K19").Value = 990949.5
$ws.Range("L19").Value = 6249.5
$ws.Range("M19").Value = -990720.5
$ws.Range("N19").Value = -6707.5
$ws.Range("H32").Value = 1172.3334
$ws.Range("I32").Value = 1241.7858
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 1241.7858
$ws.Range("L32").Value = 200
$ws.Range("M32").Value = -954.7858000000001
$ws.Range("N32").Value = -774
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H102").Value = 1420
$ws.Range("I102").Value = 1275
$ws.Range("K102").Value = 1275
$ws.Range("M102").Value = 347
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 4500
$ws.Range("K116").Value = 4500
$ws.Range("M116").Value = -2206
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("N122").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = 4500
$ws.Range("K3").Value = 4500
$ws.Range("M3").Value = -4386
$ws.Range("H4").Value = 27.5
$ws.Range("I4").Value = 25
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 25
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 90
$ws.Range("N4").Value = -260
$ws.Range("H7").Value = 1666983.4
$ws.Range("I7").Value = 2500250
$ws.Range("K7").Value = 2500250
$ws.Range("M7").Value = -2500137
$ws.Range("H134").Value = 5273.5835
$ws.Range("I134").Value = 2285.375
$ws.Range("K134").Value = 6856.125
$ws.Range("M134").Value = -4321.125

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 45.81818
$ws.Range("I7").Value = 25
$ws.Range("J7").Value = 70.8
$ws.Range("K7").Value = 25
$ws.Range("L7").Value = 70.8
$ws.Range("M7").Value = 88
$ws.Range("N7").Value = -296.8
$ws.Range("H100").Value = 99999
$ws.Range("J100").Value = 99999
$ws.Range("L100").Value = 99999
$ws.Range("N100").Value = -102163
$ws.Range("H105").Value = 7669.6665
$ws.Range("I105").Value = 14998
$ws.Range("J105").Value = 4005.5
$ws.Range("K105").Value = 14998
$ws.Range("L105").Value = 4005.5
$ws.Range("M105").Value = -13251
$ws.Range("N105").Value = -7499.5
$ws.Range("H122").Value = 2510
$ws.Range("I122").Value = 2076.75
$ws.Range("J122").Value = 3087.6667
$ws.Range("K122").Value = 6230.25
$ws.Range("L122").Value = 9263.000100000001
$ws.Range("M122").Value = -3780.25
$ws.Range("N122").Value = -14163.0001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6338
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6204

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2713.375
$ws.Range("I102").Value = 2492.9167
$ws.Range("J102").Value = 3374.75
$ws.Range("K102").Value = 2492.9167
$ws.Range("L102").Value = 3374.75
$ws.Range("M102").Value = -870.9167000000002
$ws.Range("N102").Value = -6618.75
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H122").Value = 2047.2
$ws.Range("I122").Value = 1641.75
$ws.Range("K122").Value = 4925.25
$ws.Range("M122").Value = -2475.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H10").Value = 3852
$ws.Range("I10").Value = 3800
$ws.Range("K10").Value = 3800
$ws.Range("M10").Value = -3660
$ws.Range("H15").Value = 10000000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H19").Value = 816.6667
$ws.Range("I19").Value = 200
$ws.Range("J19").Value = 1125
$ws.Range("K19").Value = 200
$ws.Range("L19").Value = 1125
$ws.Range("M19").Value = -30
$ws.Range("N19").Value = -1465
$ws.Range("H35").Value = 3186.75
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H46").Value = 3825
$ws.Range("J46").Value = 3500
$ws.Range("L46").Value = 3500
$ws.Range("N46").Value = -3876
$ws.Range("H61").Value = 9249
$ws.Range("I61").Value = 7332.6665
$ws.Range("K61").Value = 7332.6665
$ws.Range("M61").Value = -7130.6665
$ws.Range("H113").Value = 9249
$ws.Range("I113").Value = 7332.6665
$ws.Range("K113").Value = 7332.6665
$ws.Range("M113").Value = -5162.6665
$ws.Range("H136").Value = 11125
$ws.Range("I136").Value = 7000.2
$ws.Range("J136").Value = 17999.666
$ws.Range("K136").Value = 21000.6
$ws.Range("L136").Value = 53998.99800000001
$ws.Range("M136").Value = -18450.6
$ws.Range("N136").Value = -59098.99800000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 2357
$ws.Range("I122").Value = 1750
$ws.Range("K122").Value = 5250
$ws.Range("M122").Value = -2800

Write-Host "Applied all changes"